# ---- Create / position the new worksheet ----
$wb = $excel.ActiveWorkbook
$sheet1 = $wb.Worksheets.Item(1)
$ws = $wb.Worksheets.Add([Type]::Missing, $sheet1)
$ws.Name = "Criterios de Aceptación"

# ---- Title row ----
$ws.Range("A2:B2").Merge()
$ws.Range("A2").Value = "Criterios de Aceptación"
$ws.Range("A2:B2").Font.Bold = $true
$ws.Range("A2:B2").HorizontalAlignment = -4108
$ws.Range("A2:B2").VerticalAlignment = -4108
$ws.Rows.Item(2).RowHeight = 13.8

# ---- Decorative blank cell (underlined) ----
$ws.Range("D3").Font.Underline = $true
$ws.Rows.Item(3).RowHeight = 13.8

# ---- Header row ----
$ws.Range("A4").Value = "ID HISTORIA"
$ws.Range("B4").Value = "ID"
$ws.Range("C4").Value = "TITULO"
$ws.Range("D4").Value = "DADO"
$ws.Range("E4").Value = "CUANDO"
$ws.Range("F4").Value = "ENTONCES"
$ws.Range("A4:F4").Font.Bold = $true
$ws.Rows.Item(4).RowHeight = 13.8

# ---- Data rows ----
# Row 5
$ws.Range("A5").Value = 5
$ws.Range("B5").Value = 1
$ws.Range("C5").Value = "Consulta inicial"
$ws.Range("D5").Value = "Que existan usuarios registrados en el sistema"
$ws.Range("E5").Value = "El administrador abra la ventana de Listado de Usuarios"
$ws.Range("F5").Value = "Aparecerá un listado de todos los usuarios registrados "
$ws.Rows.Item(5).RowHeight = 24

# Row 6
$ws.Range("B6").Value = 2
$ws.Range("C6").Value = "Búsqueda"
$ws.Range("D6").Value = "Que el usuario elija un campo de búsqueda de la lista desplegable y escriba un dato en el campo la caja de texto"
$ws.Range("E6").Value = "El administrador presione enter en la caja de texto o de clic al botón de buscar"
$ws.Range("F6").Value = "Aparecerá un listado de los usuarios que en el campo seleccionado contenga(al principio, al fin o en medio) el dato escrito en la caja de texto "
$ws.Rows.Item(6).RowHeight = 35.25

# Row 7
$ws.Range("B7").Value = 3
$ws.Range("C7").Value = "Búsqueda con dato vacío"
$ws.Range("D7").Value = "Que el usuario elija un campo de búsqueda de la lista desplegable y no escriba un dato en el campo la caja de texto"
$ws.Range("E7").Value = "El administrador presione enter en la caja de texto o de clic al botón de buscar"
$ws.Range("F7").Value = "Aparecerá un listado de todos los usuarios registrados "
$ws.Rows.Item(7).RowHeight = 35.05

# Row 8
$ws.Range("A8").Value = 7
$ws.Range("B8").Value = 1
$ws.Range("C8").Value = "Configuración inicial"
$ws.Range("D8").Value = "Que no exista el archivo de configuración en la ruta del ejecutable"
$ws.Range("E8").Value = "Un usuario abra el sistema"
$ws.Range("F8").Value = "Aparecerá la ventana de configuración de base de datos"
$ws.Rows.Item(8).RowHeight = 24

# Row 9
$ws.Range("B9").Value = 2
$ws.Range("C9").Value = "Validación de datos"
$ws.Range("D9").Value = "Que este en blanco uno o mas campos del formulario"
$ws.Range("E9").Value = "Un usuario presione el botón guardar"
$ws.Range("F9").Value = "Aparecerá un mensaje avisando del primer dato faltante"
$ws.Rows.Item(9).RowHeight = 23.85

# Row 10
$ws.Range("B10").Value = 3
$ws.Range("C10").Value = "Conexión incorrecta"
$ws.Range("D10").Value = "Que todos los datos estén llenos pero no sean correctos(para hacer una conexión valida)"
$ws.Range("E10").Value = "Un usuario presione el botón guardar"
$ws.Range("F10").Value = "Aparecerá un mensaje indicando la conexión incorrecta"
$ws.Rows.Item(10).RowHeight = 24

# Row 11
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = "Conexión correcta"
$ws.Range("D11").Value = "Que todos los datos estén llenos y sean correctos(para hacer una conexión valida)"
$ws.Range("E11").Value = "Un usuario presione el botón guardar"
$ws.Range("F11").Value = "Aparecerá un mensaje indicando la conexión correcta y se guardaran los datos en el archivo de configuración en la ruta del ejecutable"
$ws.Rows.Item(11).RowHeight = 35.25

# Row 12
$ws.Range("B12").Value = 5
$ws.Range("C12").Value = "Cancelación de configuración"
$ws.Range("D12").Value = "Que no exista el archivo de configuración en la ruta del ejecutable"
$ws.Range("E12").Value = "El usuario presiona el botón cancelar"
$ws.Range("F12").Value = "Se cerrará el sistema"
$ws.Rows.Item(12).RowHeight = 23.85

# Row 13
$ws.Range("B13").Value = 6
$ws.Range("C13").Value = "Edición de la configuración"
$ws.Range("D13").Value = "Un administrador ingresado dentro del sistema"
$ws.Range("E13").Value = "El administrador abra la ventana de Configuración de Base de Datos"
$ws.Range("F13").Value = "Aparecerá la ventana de configuración de base de datos y cargará los datos del archivo de configuración a excepción del password"
$ws.Rows.Item(13).RowHeight = 35.25

# Row 14
$ws.Range("B14").Value = 7
$ws.Range("C14").Value = "Guardado de configuración editada"
$ws.Range("D14").Value = "Que el administrador haya escrito los datos de una conexión valida"
$ws.Range("E14").Value = "El administrador presione el botón guardar"
$ws.Range("F14").Value = "Aparecerá un mensaje indicando la conexión correcta, se guardaran los datos en el archivo de configuración y aparecerá un mensaje de reinicio del sistema, luego se reiniciará el sistema"
$ws.Rows.Item(14).RowHeight = 46.5

# Row 15
$ws.Range("B15").Value = 8
$ws.Range("C15").Value = "Prueba de conexión incorrecta"
$ws.Range("D15").Value = "Que todos los datos estén llenos y sean incorrectos(para hacer una conexión valida)"
$ws.Range("E15").Value = "El administrador presione el botón probar"
$ws.Range("F15").Value = "Aparecerá un mensaje de conexión incorrecta"
$ws.Rows.Item(15).RowHeight = 35.25

# Row 16
$ws.Range("B16").Value = 9
$ws.Range("C16").Value = "Prueba de conexión correcta"
$ws.Range("D16").Value = "Que todos los datos estén llenos y sean correctos(para hacer una conexión valida)"
$ws.Range("E16").Value = "El administrador presione el botón probar"
$ws.Range("F16").Value = "Aparecerá un mensaje de conexión correcta"
$ws.Rows.Item(16).RowHeight = 24

# ---- Merge the ID HISTORIA column for each user-story block ----
$ws.Range("A5:A7").Merge()
$ws.Range("A8:A16").Merge()

# ---- Alignment / wrap formatting ----
foreach ($a in @("A4","B4","C4","D4","E4","F4")) {
  $ws.Range($a).VerticalAlignment = -4107
}
foreach ($a in @("A5","A8")) {
  $ws.Range($a).VerticalAlignment = -4108
  $ws.Range($a).HorizontalAlignment = -4108
}
foreach ($a in @("B5","C5","B6","C6","B7","B8","E8","B9","D9","E9","B10","E10","B11","C11","E11","B12","E12","B13","D13","B14","E14","B15","E15","B16","E16")) {
  $ws.Range($a).VerticalAlignment = -4108
}
foreach ($a in @("D5")) {
  $ws.Range($a).VerticalAlignment = -4108
  $ws.Range($a).HorizontalAlignment = -4131
  $ws.Range($a).WrapText = $true
}
foreach ($a in @("E5","F5","D6","E6","F6","C7","D7","E7","F7","C8","D8","F8","C9","F9","C10","D10","F10","D11","F11","C12","D12","F12","C13","E13","F13","C14","D14","F14","C15","D15","F15","C16","D16","F16")) {
  $ws.Range($a).VerticalAlignment = -4108
  $ws.Range($a).WrapText = $true
}

# ---- Column widths (character units, converted from target OOXML widths) ----
$ws.Columns.Item(3).ColumnWidth = 15.006666666666666
$ws.Columns.Item(4).ColumnWidth = 45.43666666666667
$ws.Columns.Item(5).ColumnWidth = 35.56666666666666
$ws.Columns.Item(6).ColumnWidth = 40.016666666666666

# ---- Sheet view / selection ----
$ws.Range("E17").Select()

# ---- Update the "Historias de Usuario" sheet selection ----
$sheet1.Range("D11").Select()

# ---- Make the new sheet the active one ----
$ws.Activate()
